$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.7955796666666667
$ws.Cells.Item(2, 8).Value = 2.386739
$ws.Cells.Item(2, 9).Value = 0.1186174580157865
$ws.Cells.Item(2, 10).Value = 0.1186174580157865
$ws.Cells.Item(2, 13).Value = 145.7007446666667
$ws.Cells.Item(2, 14).Value = 437.1022340000001
$ws.Cells.Item(2, 15).Value = 0.2865937750105843
$ws.Cells.Item(2, 16).Value = 0.2865937750105843
$ws.Cells.Item(2, 17).Value = 115.9165498749918
$ws.Cells.Item(2, 18).Value = 1043.248948874926
$ws.Cells.Item(2, 19).Value = 0.03399502507490375
$ws.Cells.Item(2, 20).Value = 0.03399502507490375
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.7955796666666667
$ws.Cells.Item(3, 8).Value = 2.386739
$ws.Cells.Item(3, 9).Value = 0.1186174580157865
$ws.Cells.Item(3, 10).Value = 0.1186174580157865
$ws.Cells.Item(3, 15).Value = 0.3320294904365841
$ws.Cells.Item(3, 16).Value = 0.3320294904365841
$ws.Cells.Item(3, 17).Value = 134.2936111809791
$ws.Cells.Item(3, 18).Value = 1208.642500628812
$ws.Cells.Item(3, 19).Value = 0.03938449414186452
$ws.Cells.Item(3, 20).Value = 0.03938449414186451
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.7955796666666667
$ws.Cells.Item(4, 8).Value = 2.386739
$ws.Cells.Item(4, 9).Value = 0.1186174580157865
$ws.Cells.Item(4, 10).Value = 0.1186174580157865
$ws.Cells.Item(4, 13).Value = 128.1261546666667
$ws.Cells.Item(4, 14).Value = 384.378464
$ws.Cells.Item(4, 15).Value = 0.2520245069956105
$ws.Cells.Item(4, 16).Value = 0.2520245069956105
$ws.Cells.Item(4, 17).Value = 101.9345634209885
$ws.Cells.Item(4, 18).Value = 917.411070788896
$ws.Cells.Item(4, 19).Value = 0.02989450637750113
$ws.Cells.Item(4, 20).Value = 0.02989450637750113
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.7955796666666667
$ws.Cells.Item(5, 8).Value = 2.386739
$ws.Cells.Item(5, 9).Value = 0.1186174580157865
$ws.Cells.Item(5, 10).Value = 0.1186174580157865
$ws.Cells.Item(5, 13).Value = 65.761079
$ws.Cells.Item(5, 14).Value = 197.283237
$ws.Cells.Item(5, 15).Value = 0.1293522275572212
$ws.Cells.Item(5, 16).Value = 0.1293522275572212
$ws.Cells.Item(5, 17).Value = 52.31817731046033
$ws.Cells.Item(5, 18).Value = 470.863595794143
$ws.Cells.Item(5, 19).Value = 0.01534343242151716
$ws.Cells.Item(5, 20).Value = 0.01534343242151716
$ws.Cells.Item(6, 9).Value = 0.6312226244877757
$ws.Cells.Item(6, 10).Value = 0.6312226244877758
$ws.Cells.Item(6, 13).Value = 145.7007446666667
$ws.Cells.Item(6, 14).Value = 437.1022340000001
$ws.Cells.Item(6, 15).Value = 0.2865937750105843
$ws.Cells.Item(6, 16).Value = 0.2865937750105843
$ws.Cells.Item(6, 17).Value = 616.8497458773948
$ws.Cells.Item(6, 18).Value = 5551.647712896553
$ws.Cells.Item(6, 19).Value = 0.1809044748240401
$ws.Cells.Item(6, 20).Value = 0.1809044748240401
$ws.Cells.Item(7, 9).Value = 0.6312226244877757
$ws.Cells.Item(7, 10).Value = 0.6312226244877758
$ws.Cells.Item(7, 15).Value = 0.3320294904365841
$ws.Cells.Item(7, 16).Value = 0.3320294904365841
$ws.Cells.Item(7, 19).Value = 0.2095845263607194
$ws.Cells.Item(7, 20).Value = 0.2095845263607195
$ws.Cells.Item(8, 9).Value = 0.6312226244877757
$ws.Cells.Item(8, 10).Value = 0.6312226244877758
$ws.Cells.Item(8, 13).Value = 128.1261546666667
$ws.Cells.Item(8, 14).Value = 384.378464
$ws.Cells.Item(8, 15).Value = 0.2520245069956105
$ws.Cells.Item(8, 16).Value = 0.2520245069956105
$ws.Cells.Item(8, 17).Value = 542.4446259845547
$ws.Cells.Item(8, 18).Value = 4882.001633860992
$ws.Cells.Item(8, 19).Value = 0.159083570741007
$ws.Cells.Item(8, 20).Value = 0.159083570741007
$ws.Cells.Item(9, 9).Value = 0.6312226244877757
$ws.Cells.Item(9, 10).Value = 0.6312226244877758
$ws.Cells.Item(9, 13).Value = 65.761079
$ws.Cells.Item(9, 14).Value = 197.283237
$ws.Cells.Item(9, 15).Value = 0.1293522275572212
$ws.Cells.Item(9, 16).Value = 0.1293522275572212
$ws.Cells.Item(9, 17).Value = 278.411101896404
$ws.Cells.Item(9, 18).Value = 2505.699917067636
$ws.Cells.Item(9, 19).Value = 0.08165005256200916
$ws.Cells.Item(9, 20).Value = 0.08165005256200918
$ws.Cells.Item(10, 7).Value = 1.666370333333333
$ws.Cells.Item(10, 8).Value = 4.999111
$ws.Cells.Item(10, 9).Value = 0.248448548064433
$ws.Cells.Item(10, 10).Value = 0.248448548064433
$ws.Cells.Item(10, 13).Value = 145.7007446666667
$ws.Cells.Item(10, 14).Value = 437.1022340000001
$ws.Cells.Item(10, 15).Value = 0.2865937750105843
$ws.Cells.Item(10, 16).Value = 0.2865937750105843
$ws.Cells.Item(10, 17).Value = 242.7913984571082
$ws.Cells.Item(10, 18).Value = 2185.122586113974
$ws.Cells.Item(10, 19).Value = 0.07120380728568443
$ws.Cells.Item(10, 20).Value = 0.07120380728568443
$ws.Cells.Item(11, 7).Value = 1.666370333333333
$ws.Cells.Item(11, 8).Value = 4.999111
$ws.Cells.Item(11, 9).Value = 0.248448548064433
$ws.Cells.Item(11, 10).Value = 0.248448548064433
$ws.Cells.Item(11, 15).Value = 0.3320294904365841
$ws.Cells.Item(11, 16).Value = 0.3320294904365841
$ws.Cells.Item(11, 17).Value = 281.2828167992209
$ws.Cells.Item(11, 18).Value = 2531.545351192988
$ws.Cells.Item(11, 19).Value = 0.08249224481354285
$ws.Cells.Item(11, 20).Value = 0.08249224481354286
$ws.Cells.Item(12, 7).Value = 1.666370333333333
$ws.Cells.Item(12, 8).Value = 4.999111
$ws.Cells.Item(12, 9).Value = 0.248448548064433
$ws.Cells.Item(12, 10).Value = 0.248448548064433
$ws.Cells.Item(12, 13).Value = 128.1261546666667
$ws.Cells.Item(12, 14).Value = 384.378464
$ws.Cells.Item(12, 15).Value = 0.2520245069956105
$ws.Cells.Item(12, 16).Value = 0.2520245069956105
$ws.Cells.Item(12, 17).Value = 213.5056230606116
$ws.Cells.Item(12, 18).Value = 1921.550607545504
$ws.Cells.Item(12, 19).Value = 0.06261512283971395
$ws.Cells.Item(12, 20).Value = 0.06261512283971395
$ws.Cells.Item(13, 7).Value = 1.666370333333333
$ws.Cells.Item(13, 8).Value = 4.999111
$ws.Cells.Item(13, 9).Value = 0.248448548064433
$ws.Cells.Item(13, 10).Value = 0.248448548064433
$ws.Cells.Item(13, 13).Value = 65.761079
$ws.Cells.Item(13, 14).Value = 197.283237
$ws.Cells.Item(13, 15).Value = 0.1293522275572212
$ws.Cells.Item(13, 16).Value = 0.1293522275572212
$ws.Cells.Item(13, 17).Value = 109.5823111335897
$ws.Cells.Item(13, 18).Value = 986.2408002023069
$ws.Cells.Item(13, 19).Value = 0.03213737312549175
$ws.Cells.Item(13, 20).Value = 0.03213737312549175
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.01147833333333333
$ws.Cells.Item(14, 8).Value = 0.034435
$ws.Cells.Item(14, 9).Value = 0.00171136943200476
$ws.Cells.Item(14, 10).Value = 0.00171136943200476
$ws.Cells.Item(14, 13).Value = 145.7007446666667
$ws.Cells.Item(14, 14).Value = 437.1022340000001
$ws.Cells.Item(14, 15).Value = 0.2865937750105843
$ws.Cells.Item(14, 16).Value = 0.2865937750105843
$ws.Cells.Item(14, 17).Value = 1.672401714198889
$ws.Cells.Item(14, 18).Value = 15.05161542779
$ws.Cells.Item(14, 19).Value = 0.0004904678259559636
$ws.Cells.Item(14, 20).Value = 0.0004904678259559636
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.01147833333333333
$ws.Cells.Item(15, 8).Value = 0.034435
$ws.Cells.Item(15, 9).Value = 0.00171136943200476
$ws.Cells.Item(15, 10).Value = 0.00171136943200476
$ws.Cells.Item(15, 15).Value = 0.3320294904365841
$ws.Cells.Item(15, 16).Value = 0.3320294904365841
$ws.Cells.Item(15, 17).Value = 1.937539253775556
$ws.Cells.Item(15, 18).Value = 17.43785328398
$ws.Cells.Item(15, 19).Value = 0.0005682251204572869
$ws.Cells.Item(15, 20).Value = 0.0005682251204572869
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.01147833333333333
$ws.Cells.Item(16, 8).Value = 0.034435
$ws.Cells.Item(16, 9).Value = 0.00171136943200476
$ws.Cells.Item(16, 10).Value = 0.00171136943200476
$ws.Cells.Item(16, 13).Value = 128.1261546666667
$ws.Cells.Item(16, 14).Value = 384.378464
$ws.Cells.Item(16, 15).Value = 0.2520245069956105
$ws.Cells.Item(16, 16).Value = 0.2520245069956105
$ws.Cells.Item(16, 17).Value = 1.470674711982222
$ws.Cells.Item(16, 18).Value = 13.23607240784
$ws.Cells.Item(16, 19).Value = 0.0004313070373883576
$ws.Cells.Item(16, 20).Value = 0.0004313070373883576
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 0.6666666666666666
$ws.Cells.Item(17, 7).Value = 0.01147833333333333
$ws.Cells.Item(17, 8).Value = 0.034435
$ws.Cells.Item(17, 9).Value = 0.00171136943200476
$ws.Cells.Item(17, 10).Value = 0.00171136943200476
$ws.Cells.Item(17, 13).Value = 65.761079
$ws.Cells.Item(17, 14).Value = 197.283237
$ws.Cells.Item(17, 15).Value = 0.1293522275572212
$ws.Cells.Item(17, 16).Value = 0.1293522275572212
$ws.Cells.Item(17, 17).Value = 0.7548275851216666
$ws.Cells.Item(17, 18).Value = 6.793448266095
$ws.Cells.Item(17, 19).Value = 0.0002213694482031522
$ws.Cells.Item(17, 20).Value = 0.0002213694482031522